$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Remove the old "X/Y voulue/actuel" scratch calculation block (K21:L33) ---
$ws.Range("K21:L33").ClearContents()

# --- Add the new "erreur distance / erreur orientation" block (G36:I51) ---
$ws.Range("G36").Value = "x cons"
$ws.Range("H36").Value = 100

$ws.Range("G37").Value = "y cons"
$ws.Range("H37").Value = 100

$ws.Range("G38").Value = "x actu"
$ws.Range("H38").Value = 0

$ws.Range("G39").Value = "y act"
$ws.Range("H39").Value = 0

$ws.Range("G41").Value = "x cons - x actu"
$ws.Range("H41").Formula = "=H36-H38"

$ws.Range("G42").Value = "y cons - y actu"
$ws.Range("H42").Formula = "=H37-H39"

$ws.Range("G44").Value = "X au carré"
$ws.Range("H44").Formula = "=H41*H41"

$ws.Range("G45").Value = "y au carré"
$ws.Range("H45").Formula = "=H42*H42"

$ws.Range("G46").Value = "Erreur dist"
$ws.Range("H46").Formula = "=SQRT(H44+H45)"

$ws.Range("G48").Value = "theta robot"
$ws.Range("H48").Value = -1.039
$ws.Range("I48").Value = "deg"

$ws.Range("H49").Formula = "=RADIANS(H48)"
$ws.Range("I49").Value = "rad"

$ws.Range("G50").Value = "erreur orient"
$ws.Range("H50").Formula = "=ATAN2(H41,H42)-H49"
$ws.Range("I50").Value = "rad"

$ws.Range("H51").Formula = "=DEGREES(H50)"
$ws.Range("I51").Value = "deg"

# --- Update the view: scrolled/selected area moved down to the new block ---
[void]$ws.Activate()
[void]$ws.Range("G47").Select()
$win = $excel.ActiveWindow()
$win.ScrollRow = 41
[void]($win.ScrollColumn = 7)
